$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data (session 12, status 60)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 60

# Update selection to B14, matching the post-entry cursor position
$ws.Range("B14").Select()
